$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.527.14'
$ws.Range("E2").Value = '  -2.72%  '
$ws.Range("D3").Value = '1.814.52'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("D4").Value = '''1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '''308.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '''1.008'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("D7").Value = '''0.4565'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").Value = '''0.3670'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").Value = '''0.07140'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("D10").Value = '''0.8815'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '''0.07757'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '''19.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.70%  '
$ws.Range("D13").Value = '1.835.02'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").Value = '''5.300'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").Value = '''6.376'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '''86.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.34%  '
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = '''0.000008605'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.68%  '
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = '26.586.27'
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("D21").Value = '''14.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.21%  '
$ws.Range("D22").Value = '''5.021'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("D23").Value = '''10.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = '''1.985'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").Value = '''151.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").Value = '''17.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("D27").Value = '''2.076'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '''113.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("D29").Value = '''4.865'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.61%  '
$ws.Range("D30").Value = '''0.08698'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.69%  '
$ws.Range("D31").Value = '''3.040'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("D32").Value = '''4.504'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").Value = '''0.7333'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.93%  '
$ws.Range("D34").Value = '''1.121'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.20%  '
$ws.Range("D35").Value = '''2.680'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '''1.086'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("D38").Value = '''0.01961'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.88%  '
$ws.Range("D39").Value = '''0.05133'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = '''2.893'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.77%  '
$ws.Range("D41").Value = '''7.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = '''0.5007'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").Value = '''0.1557'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.44%  '
$ws.Range("D44").Value = '''8.164'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.24%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = '''0.4607'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.02%  '
$ws.Range("D47").Value = '''9.989'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("D48").Value = '''101.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("D49").Value = '''1.591'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").Value = '''0.06005'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").Value = '''64.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.57%  '
